# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2183
    $ws.Range("F4").Value = 1561
    $ws.Range("F5").Value = 7323
    $ws.Range("F7").Value = 176
}
